# Integration plan update: combine steps 9-11 into a single step 9.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# New step 9 (row 11) is now fully complete ("X" in every column),
# absorbing what used to be steps 10 and 11.
$ws.Range("H11").Value = "X"
$ws.Range("I11").Value = "X"

# Old rows for steps 10 and 11 are no longer needed - remove them entirely.
$ws.Rows("12:13").Delete()
